# Marksheet update: fill in "Marking" row as a real number (not text), recompute the
# scoring summary (Right/Wrong/Max + totals), and populate the "Student Ans" column
# for the quiz's first answer block (col A, and the first three rows of the second
# block in col D) with colour-coded correct/incorrect answers. The now-unused
# second/third "Student Ans" helper columns (D19:E40 and the whole G:H block) are
# removed, which also shrinks the sheet's used range from A5:H40 down to A5:E40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlHAlignCenter = -4108

# ---------------------------------------------------------------------------
# Summary block (rows 10-12): Right / Wrong / Not Attempt / Max, and totals.
# ---------------------------------------------------------------------------

# Row/Col labels ("No.", "Marking", "Total") pick up the same centred title
# style used elsewhere in the sheet.
foreach ($r in 10, 11, 12) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Style = "mtitleStyle"
    $cell.HorizontalAlignment = $xlHAlignCenter
}

# Right / Wrong / Not Attempt / Max counts for this answer block.
$ws.Range("B10").Value = 12
$ws.Range("C10").Value = 16
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 28

# Marking scheme: +4 for right, -1 for wrong. The wrong-answer mark used to be
# typed in as the text "-1"; store it as a real (float-capable) number instead.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Totals: 12*4 = 48, 16*-1 = -16, giving a final score of 32 out of 28*4 = 112.
$ws.Range("B12").Value = 48
$ws.Range("C12").Value = -16
$ws.Range("E12").Value = "32/112"

# ---------------------------------------------------------------------------
# Student answers for the first answer block (column A) plus the leftover
# first three rows of the second block (column D), colour-coded against the
# "Correct Ans" column next to each (B for col A, E for col D).
# ---------------------------------------------------------------------------

$answers = @(
  @{Cell="A16"; Val="Option A"; Style="correctStyle"},
  @{Cell="A17"; Val="Option D"; Style="correctStyle"},
  @{Cell="A18"; Val="Option C"; Style="incorrectStyle"},
  @{Cell="A19"; Val="Option C"; Style="correctStyle"},
  @{Cell="A20"; Val="Option D"; Style="incorrectStyle"},
  @{Cell="A21"; Val="Option A"; Style="incorrectStyle"},
  @{Cell="A22"; Val="Option D"; Style="correctStyle"},
  @{Cell="A23"; Val="Option A"; Style="incorrectStyle"},
  @{Cell="A24"; Val="Option C"; Style="incorrectStyle"},
  @{Cell="A25"; Val="Option A"; Style="correctStyle"},
  @{Cell="A26"; Val="Option B"; Style="incorrectStyle"},
  @{Cell="A27"; Val="Option A"; Style="correctStyle"},
  @{Cell="A28"; Val="Option C"; Style="incorrectStyle"},
  @{Cell="A29"; Val="Option B"; Style="incorrectStyle"},
  @{Cell="A30"; Val="Option B"; Style="correctStyle"},
  @{Cell="A31"; Val="Option B"; Style="incorrectStyle"},
  @{Cell="A32"; Val="Option C"; Style="correctStyle"},
  @{Cell="A33"; Val="Option A"; Style="incorrectStyle"},
  @{Cell="A34"; Val="Option D"; Style="incorrectStyle"},
  @{Cell="A35"; Val="Option D"; Style="correctStyle"},
  @{Cell="A36"; Val="Option B"; Style="incorrectStyle"},
  @{Cell="A37"; Val="Option A"; Style="correctStyle"},
  @{Cell="A38"; Val="Option B"; Style="incorrectStyle"},
  @{Cell="A39"; Val="Option D"; Style="correctStyle"},
  @{Cell="A40"; Val="Option B"; Style="incorrectStyle"},
  @{Cell="D16"; Val="Option C"; Style="incorrectStyle"},
  @{Cell="D17"; Val="Option C"; Style="correctStyle"},
  @{Cell="D18"; Val="Option B"; Style="incorrectStyle"}
)

foreach ($item in $answers) {
    $cell = $ws.Range($item.Cell)
    $cell.Value = $item.Val
    $cell.Style = $item.Style
    $cell.HorizontalAlignment = $xlHAlignCenter
}

# ---------------------------------------------------------------------------
# Drop the now-unused helper columns: the whole third "Student Ans"/"Correct
# Ans" block (G:H, rows 15-21) and the remainder of the second block
# (D19:E40) that no longer has matching student answers.
# ---------------------------------------------------------------------------

$ws.Range("G15:H21").Clear()
$ws.Range("D19:E40").Clear()
